$wb = $excel.ActiveWorkbook

# --- Sheet: P_valores ---
$ws = $wb.Worksheets.Item("P_valores")

$ws.Range("C2").Value = 0.5381240820189352
$ws.Range("D2").Value = 0.5082764740909798
$ws.Range("E2").Value = 0.9609948908993071
$ws.Range("F2").Value = 0.1621277125216749

$ws.Range("B3").Value = 0.5381240820189352
$ws.Range("D3").Value = 0.1039655594518818
$ws.Range("E3").Value = 0.5719463725390876
$ws.Range("F3").Value = 0.267637115530216

$ws.Range("B4").Value = 0.5082764740909798
$ws.Range("C4").Value = 0.1039655594518818
$ws.Range("E4").Value = 0.3058389705629412
$ws.Range("F4").Value = 0.01710459219489247

$ws.Range("B5").Value = 0.9609948908993071
$ws.Range("C5").Value = 0.5719463725390876
$ws.Range("D5").Value = 0.3058389705629412
$ws.Range("F5").Value = 0.2330005773306945

$ws.Range("B6").Value = 0.1621277125216749
$ws.Range("C6").Value = 0.267637115530216
$ws.Range("D6").Value = 0.01710459219489247
$ws.Range("E6").Value = 0.2330005773306945

# --- Sheet: Estadisticos_DM ---
$ws2 = $wb.Worksheets.Item("Estadisticos_DM")

$ws2.Range("C2").Value = -0.6311295421349024
$ws2.Range("D2").Value = 0.6788889858102605
$ws2.Range("E2").Value = -0.04978773092889133
$ws2.Range("F2").Value = -1.475818022211977

$ws2.Range("B3").Value = 0.6311295421349024
$ws2.Range("D3").Value = 1.738997072023389
$ws2.Range("E3").Value = 0.5787586386040297
$ws2.Range("F3").Value = -1.154462259052182

$ws2.Range("B4").Value = -0.6788889858102605
$ws2.Range("C4").Value = -1.738997072023389
$ws2.Range("E4").Value = -1.062853138435228
$ws2.Range("F4").Value = -2.704518992468954

$ws2.Range("B5").Value = 0.04978773092889133
$ws2.Range("C5").Value = -0.5787586386040297
$ws2.Range("D5").Value = 1.062853138435228
$ws2.Range("F5").Value = -1.24659585269093

$ws2.Range("B6").Value = 1.475818022211977
$ws2.Range("C6").Value = 1.154462259052182
$ws2.Range("D6").Value = 2.704518992468954
$ws2.Range("E6").Value = 1.24659585269093
